# file düzenlendi. import from excel düzenlendi
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet (KAYIT). This also repoints the _FilterDatabase defined name's sheet part.
$ws.Name = "KAYIT"

# 2) Copy row 2's cell formatting down onto rows 3:8 (reuses existing style indices 4/5/6).
$ws.Range("A2:I2").Copy()
$ws.Range("A3:I8").PasteSpecial(-4122)

# 3) Column width tweaks for F, G, H.
$ws.Columns.Item(6).ColumnWidth = 46
$ws.Columns.Item(7).ColumnWidth = 46.8333333333333
$ws.Columns.Item(8).ColumnWidth = 43.8333333333333

# 4) Row heights: row2 grows from 35.25 to 51; rows 3:8 also become 51.
$ws.Range("A2:A8").RowHeight = 51

# 5) Fill in the three new distributor rows (3, 4, 5).
#    A/B/C mirror row 2 (KUZEY / ASYA / 100); E holds the distributor code.
$ws.Cells.Item(3,1).Value = "KUZEY"
$ws.Cells.Item(3,2).Value = "ASYA"
$ws.Cells.Item(3,3).Value = 100
$ws.Cells.Item(3,5).Value = 605

$ws.Cells.Item(4,1).Value = "KUZEY"
$ws.Cells.Item(4,2).Value = "ASYA"
$ws.Cells.Item(4,3).Value = 100
$ws.Cells.Item(4,5).Value = 606

$ws.Cells.Item(5,1).Value = "KUZEY"
$ws.Cells.Item(5,2).Value = "ASYA"
$ws.Cells.Item(5,3).Value = 100
$ws.Cells.Item(5,5).Value = 607

# Distributor / contact text, entered row-by-row (4, 5, then 3) in D, H, G, F order -
# matches the order the new shared strings were authored in the workbook.
$ws.Cells.Item(4,4).Value = "KOPUZ DIŞ"
$ws.Cells.Item(4,8).Value = "Nuray ÇELİK<nuraycelik@kopuz.com.tr>"
$ws.Cells.Item(4,7).Value = "Taner MARANGOZ<tanermarangoz@kopuz.com.tr>,
Cemal KOPUZ<cemal@kopuz.com.tr>,
Ahmet USTAOĞLU<ahmet.ustaoglu@kopuz.com.tr>"
$ws.Cells.Item(4,6).Value = "Özgür YILDIRIM<ozguryildirim@kopuz.com.tr>,
Yavuz ELMAS<yavuzelmas@kopuz.com.tr>,
Yalçın SÜZMETAŞ<yalcinsuzmetas@kopuz.com.tr>"

$ws.Cells.Item(5,4).Value = "KOPUZ İÇ"
$ws.Cells.Item(5,8).Value = "Seval TUNCER<seval.tuncer@kopuz.com.tr>"
$ws.Cells.Item(5,7).Value = "Tuncay KARAKAYA<tuncaykarakaya@kopuz.com.tr>,
Vedat YURTSEVEN<vedatyurduseven@kopuz.com.tr>"
$ws.Cells.Item(5,6).Value = "Emre ÇİFTÇİOĞLU<emreciftcioglu@kopuz.com.tr>,
İsmail KILIÇASLAN<ismailkilicaslan@kopuz.com.tr>"

$ws.Cells.Item(3,4).Value = "ASYA PAZARLAMA"
$ws.Cells.Item(3,8).Value = "Ramazan SARIYILDIZ<asyapazarlama@gmail.com>"
$ws.Cells.Item(3,6).Value = "Melih BIYIKLI<melihbiyikli@gmail.com>"
# G3 (TTE) intentionally left blank.

# 6) Append an extra blank row (12) with the same style as row 11.
$ws.Range("A11:I11").Copy()
$ws.Range("A12:I12").PasteSpecial(-4122)

# 7) Fix up the hidden _FilterDatabase defined name (row count 241 -> 242).
foreach ($n in $wb.Names) {
  if ($n.Name -eq "KAYIT!_FilterDatabase") {
    $n.RefersTo = "=KAYIT!`$A`$1:`$H`$242"
  }
}

# 8) Leave the selection on F4, matching the saved view state.
$ws.Range("F4").Select()
